# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the newly generated output, per the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 5604
$wsExhibit.Range("F5").Value = 684
$wsExhibit.Range("F6").Value = 676
$wsExhibit.Range("F7").Value = 30
$wsExhibit.Range("F12").Value = 5441
$wsExhibit.Range("F13").Value = 465
$wsExhibit.Range("F14").Value = 291
$wsExhibit.Range("F15").Value = 252
$wsExhibit.Range("F16").Value = 58
$wsExhibit.Range("F17").Value = 26
$wsExhibit.Range("F19").Value = 4585
$wsExhibit.Range("F20").Value = 234
$wsExhibit.Range("F21").Value = 1215
$wsExhibit.Range("F26").Value = 206
$wsExhibit.Range("F37").Value = 46
$wsExhibit.Range("F38").Value = 52

# --- Sheet "全部类型" (sheet4) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5604
$wsAll.Range("F5").Value = 684
$wsAll.Range("F6").Value = 676
$wsAll.Range("F7").Value = 30
$wsAll.Range("F11").Value = 1589
$wsAll.Range("F12").Value = 5441
$wsAll.Range("F13").Value = 465
$wsAll.Range("F14").Value = 291
$wsAll.Range("F15").Value = 252
$wsAll.Range("F16").Value = 58
$wsAll.Range("F17").Value = 26
$wsAll.Range("F19").Value = 4585
$wsAll.Range("F20").Value = 234
$wsAll.Range("F21").Value = 1215
$wsAll.Range("F26").Value = 206
$wsAll.Range("F37").Value = 46
$wsAll.Range("F38").Value = 52

$wb.Save()
